# Auto-generated Excel COM-interop script applying numeric updates
# to the Jenova_Profits workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 96 (ALC)
$ws.Range("H96").Value = 732.8889
$ws.Range("I96").Value = 656.7143
$ws.Range("K96").Value = 1970.1429
$ws.Range("M96").Value = -597.1428999999998

# Row 103 (ALC)
$ws.Range("H103").Value = 800.5
$ws.Range("I103").Value = 2004
$ws.Range("K103").Value = 6012
$ws.Range("M103").Value = -5426

# Row 112 (ALC)
$ws.Range("H112").Value = 3238.0952
$ws.Range("J112").Value = 3238.0952
$ws.Range("L112").Value = 9714.285600000001
$ws.Range("N112").Value = -11930.2856

# Row 125 (ALC)
$ws.Range("H125").Value = 9262644
$ws.Range("J125").Value = 11114801
$ws.Range("L125").Value = 100033209
$ws.Range("N125").Value = -100038129

# Row 132 (ALC)
$ws.Range("H132").Value = 2514.5483
$ws.Range("I132").Value = 2377.6206
$ws.Range("K132").Value = 7132.861800000001
$ws.Range("M132").Value = -4602.861800000001

# Row 137 (ALC)
$ws.Range("H137").Value = 6818.9375
$ws.Range("I137").Value = 7278.385
$ws.Range("J137").Value = 4828
$ws.Range("K137").Value = 21835.155
$ws.Range("L137").Value = 14484
$ws.Range("M137").Value = -19285.155
$ws.Range("N137").Value = -19584

$ws = $wb.Worksheets.Item("ARM")
# Row 31 (ARM)
$ws.Range("H31").Value = 5468.2
$ws.Range("I31").Value = 5468.2
$ws.Range("K31").Value = 5468.2
$ws.Range("M31").Value = -5174.2

# Row 32 (ARM)
$ws.Range("H32").Value = 3533.0874
$ws.Range("I32").Value = 3290.1794
$ws.Range("K32").Value = 3290.1794
$ws.Range("M32").Value = -3003.1794

# Row 45 (ARM)
$ws.Range("H45").Value = 2131.6155
$ws.Range("I45").Value = 1445.6666
$ws.Range("K45").Value = 1445.6666
$ws.Range("M45").Value = -1068.6666

# Row 61 (ARM)
$ws.Range("H61").Value = 4732.467
$ws.Range("I61").Value = 4460.769
$ws.Range("J61").Value = 6498.5
$ws.Range("K61").Value = 4460.769
$ws.Range("L61").Value = 6498.5
$ws.Range("M61").Value = -4248.769
$ws.Range("N61").Value = -6922.5

# Row 74 (ARM)
$ws.Range("H74").Value = 2503.25
$ws.Range("I74").Value = 1992.5
$ws.Range("J74").Value = 3014
$ws.Range("K74").Value = 1992.5
$ws.Range("L74").Value = 3014
$ws.Range("M74").Value = -1118.5
$ws.Range("N74").Value = -4762

# Row 77 (ARM)
$ws.Range("H77").Value = 2503.25
$ws.Range("I77").Value = 1992.5
$ws.Range("J77").Value = 3014
$ws.Range("K77").Value = 9962.5
$ws.Range("L77").Value = 15070
$ws.Range("M77").Value = -5594.5
$ws.Range("N77").Value = -23806

# Row 102 (ARM)
$ws.Range("H102").Value = 2358.6428
$ws.Range("I102").Value = 1902.1
$ws.Range("K102").Value = 1902.1
$ws.Range("M102").Value = -280.0999999999999

# Row 132 (ARM)
$ws.Range("H132").Value = 3896.975
$ws.Range("I132").Value = 4182.731
$ws.Range("J132").Value = 3366.2856
$ws.Range("K132").Value = 12548.193
$ws.Range("L132").Value = 10098.8568
$ws.Range("M132").Value = -10018.193
$ws.Range("N132").Value = -15158.8568

# Row 136 (ARM)
$ws.Range("H136").Value = 4732.467
$ws.Range("I136").Value = 4460.769
$ws.Range("J136").Value = 6498.5
$ws.Range("K136").Value = 13382.307
$ws.Range("L136").Value = 19495.5
$ws.Range("M136").Value = -10832.307
$ws.Range("N136").Value = -24595.5

$ws = $wb.Worksheets.Item("BSM")
# Row 20 (BSM)
$ws.Range("H20").Value = 2788.7932
$ws.Range("I20").Value = 2618.7058
$ws.Range("J20").Value = 3029.75
$ws.Range("K20").Value = 2618.7058
$ws.Range("L20").Value = 3029.75
$ws.Range("M20").Value = -2371.7058
$ws.Range("N20").Value = -3523.75

# Row 22 (BSM)
$ws.Range("H22").Value = 195
$ws.Range("I22").Value = 195
$ws.Range("K22").Value = 195
$ws.Range("M22").Value = -22

# Row 86 (BSM)
$ws.Range("H86").Value = 1420746.9
$ws.Range("I86").Value = 1892777.4
$ws.Range("J86").Value = 4655.3335
$ws.Range("K86").Value = 1892777.4
$ws.Range("L86").Value = 4655.3335
$ws.Range("M86").Value = -1891654.4
$ws.Range("N86").Value = -6901.3335

# Row 89 (BSM)
$ws.Range("H89").Value = 1420746.9
$ws.Range("I89").Value = 1892777.4
$ws.Range("J89").Value = 4655.3335
$ws.Range("K89").Value = 9463887
$ws.Range("L89").Value = 23276.6675
$ws.Range("M89").Value = -9458271
$ws.Range("N89").Value = -34508.6675

# Row 134 (BSM)
$ws.Range("H134").Value = 55350.285
$ws.Range("I134").Value = 8097.4375
$ws.Range("J134").Value = 206559.4
$ws.Range("K134").Value = 24292.3125
$ws.Range("L134").Value = 619678.2
$ws.Range("M134").Value = -21757.3125
$ws.Range("N134").Value = -624748.2

$ws = $wb.Worksheets.Item("CRP")
# Row 7 (CRP)
$ws.Range("H7").Value = 55.166668
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

# Row 31 (CRP)
$ws.Range("H31").Value = 36695.332
$ws.Range("I31").Value = 1326.6
$ws.Range("K31").Value = 1326.6
$ws.Range("M31").Value = -1031.6

# Row 34 (CRP)
$ws.Range("H34").Value = 36695.332
$ws.Range("I34").Value = 1326.6
$ws.Range("K34").Value = 1326.6
$ws.Range("M34").Value = -1124.6

# Row 58 (CRP)
$ws.Range("H58").Value = 6545.5
$ws.Range("I58").Value = 6545.5
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 6545.5
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -6342.5
$ws.Range("N58").ClearContents()

# Row 94 (CRP)
$ws.Range("H94").Value = 1051.3636
$ws.Range("J94").Value = 1175.8889
$ws.Range("L94").Value = 1175.8889
$ws.Range("N94").Value = -2077.8889

# Row 132 (CRP)
$ws.Range("H132").Value = 2356.7222
$ws.Range("J132").Value = 2727.2856
$ws.Range("L132").Value = 8181.8568
$ws.Range("N132").Value = -13241.8568

# Row 133 (CRP)
$ws.Range("H133").Value = 51428.145
$ws.Range("J133").Value = 55999.6
$ws.Range("L133").Value = 55999.6
$ws.Range("N133").Value = -61059.6

# Row 134 (CRP)
$ws.Range("H134").Value = 419555.38
$ws.Range("I134").Value = 3014.3044
$ws.Range("K134").Value = 9042.913199999999
$ws.Range("M134").Value = -6507.913199999999

# Row 136 (CRP)
$ws.Range("H136").Value = 6545.5
$ws.Range("I136").Value = 6545.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 19636.5
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -17086.5
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 117 (CUL)
$ws.Range("H117").Value = 937.44446
$ws.Range("J117").Value = 1272
$ws.Range("L117").Value = 3816
$ws.Range("N117").Value = -10700

# Row 122 (CUL)
$ws.Range("H122").Value = 32451.281
$ws.Range("I122").Value = 778.93335
$ws.Range("J122").Value = 60397.47
$ws.Range("K122").Value = 7010.40015
$ws.Range("L122").Value = 543577.23
$ws.Range("M122").Value = -4560.40015
$ws.Range("N122").Value = -548477.23

# Row 131 (CUL)
$ws.Range("H131").Value = 4980.2173
$ws.Range("I131").Value = 1520.8182
$ws.Range("J131").Value = 8151.3335
$ws.Range("K131").Value = 4562.4546
$ws.Range("L131").Value = 24454.0005
$ws.Range("M131").Value = 477.5454
$ws.Range("N131").Value = -34534.00049999999

# Row 132 (CUL)
$ws.Range("H132").Value = 462805.72
$ws.Range("I132").Value = 112149.22
$ws.Range("J132").Value = 673199.6
$ws.Range("K132").Value = 1009342.98
$ws.Range("L132").Value = 6058796.399999999
$ws.Range("M132").Value = -1006812.98
$ws.Range("N132").Value = -6063856.399999999

# Row 133 (CUL)
$ws.Range("H133").Value = 7181.6924
$ws.Range("I133").Value = 5925.8
$ws.Range("J133").Value = 7966.625
$ws.Range("K133").Value = 17777.4
$ws.Range("L133").Value = 23899.875
$ws.Range("M133").Value = -12717.4
$ws.Range("N133").Value = -34019.875

$ws = $wb.Worksheets.Item("GSM")
# Row 11 (GSM)
$ws.Range("H11").Value = 12500000
$ws.Range("I11").Value = 15000000
$ws.Range("K11").Value = 15000000
$ws.Range("M11").Value = -14999861

# Row 19 (GSM)
$ws.Range("H19").Value = 99000
$ws.Range("J19").Value = 99000
$ws.Range("L19").Value = 99000
$ws.Range("N19").Value = -99576

# Row 70 (GSM)
$ws.Range("H70").Value = 10446
$ws.Range("I70").Value = 7865.4443
$ws.Range("K70").Value = 7865.4443
$ws.Range("M70").Value = -7595.4443

# Row 73 (GSM)
$ws.Range("H73").Value = 10446
$ws.Range("I73").Value = 7865.4443
$ws.Range("K73").Value = 7865.4443
$ws.Range("M73").Value = -6929.4443

# Row 93 (GSM)
$ws.Range("H93").Value = 39956.332
$ws.Range("I93").Value = 39949
$ws.Range("J93").Value = 39960
$ws.Range("K93").Value = 39949
$ws.Range("L93").Value = 39960
$ws.Range("M93").Value = -38077
$ws.Range("N93").Value = -43704

# Row 122 (GSM)
$ws.Range("H122").Value = 5666.6665
$ws.Range("I122").Value = 4400
$ws.Range("J122").Value = 6300
$ws.Range("K122").Value = 13200
$ws.Range("L122").Value = 18900
$ws.Range("M122").Value = -10750
$ws.Range("N122").Value = -23800

# Row 123 (GSM)
$ws.Range("H123").Value = 43498
$ws.Range("J123").Value = 43498
$ws.Range("L123").Value = 43498
$ws.Range("N123").Value = -48398

$ws = $wb.Worksheets.Item("LTW")
# Row 23 (LTW)
$ws.Range("H23").Value = 756125
$ws.Range("I23").Value = 756125
$ws.Range("K23").Value = 756125
$ws.Range("M23").Value = -755895

# Row 132 (LTW)
$ws.Range("H132").Value = 8453.280000000001
$ws.Range("I132").Value = 7334.6924
$ws.Range("J132").Value = 9665.083000000001
$ws.Range("K132").Value = 22004.0772
$ws.Range("L132").Value = 28995.249
$ws.Range("M132").Value = -19474.0772
$ws.Range("N132").Value = -34055.249

# Row 136 (LTW)
$ws.Range("H136").Value = 262541.28
$ws.Range("I136").Value = 480409.62
$ws.Range("K136").Value = 1441228.86
$ws.Range("M136").Value = -1438678.86

$ws = $wb.Worksheets.Item("WVR")
# Row 122 (WVR)
$ws.Range("H122").Value = 37039948
$ws.Range("J122").Value = 3965.4443
$ws.Range("L122").Value = 11896.3329
$ws.Range("N122").Value = -16796.3329

# Row 132 (WVR)
$ws.Range("H132").Value = 32840.258
$ws.Range("I132").Value = 1717.4375
$ws.Range("K132").Value = 5152.3125
$ws.Range("M132").Value = -2622.3125

# Row 136 (WVR)
$ws.Range("H136").Value = 10176177
$ws.Range("I136").Value = 12363024
$ws.Range("J136").Value = 335363.66
$ws.Range("K136").Value = 37089072
$ws.Range("L136").Value = 1006090.98
$ws.Range("M136").Value = -37086522
$ws.Range("N136").Value = -1011190.98
